$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -4
    3  = 3
    5  = -3
    6  = -3
    7  = 2
    8  = -3
    9  = 6
    11 = -2
    12 = 1
    13 = -2
    14 = -1
    15 = -3
    17 = 2
    18 = -1
    19 = -3
    20 = 4
    21 = 2
    22 = 7
    23 = -2
    24 = 3
    25 = -3
    27 = 5
    28 = -1
    29 = 7
    30 = -4
    31 = 1
    32 = 1
    33 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
